$p = $ppt.ActivePresentation

# Add a new slide at the end, using the "Title and Content" layout
# (same CustomLayout used by all the other content slides in this deck).
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# ---- Title placeholder ----
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Excel"
$title.InsertAfter("动手实验室 ") | Out-Null
$title.InsertAfter("- ") | Out-Null
$title.InsertAfter("冻") | Out-Null
$title.InsertAfter("结或锁定窗格") | Out-Null

# ---- Body / content placeholder ----
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "冻结窗格"
$body.InsertAfter("`r冻结") | Out-Null
$body.InsertAfter("首") | Out-Null
$body.InsertAfter("行") | Out-Null
$body.InsertAfter("`r冻") | Out-Null
$body.InsertAfter("结首列") | Out-Null
$body.InsertAfter("`r冻结自定义行、列") | Out-Null
$body.InsertAfter("`r打") | Out-Null
$body.InsertAfter("印") | Out-Null
$body.InsertAfter("`r每一") | Out-Null
$body.InsertAfter("页打印标题") | Out-Null
$body.InsertAfter("`r.") | Out-Null
$body.Paragraphs(7, 1).Text = ""
$body.Paragraphs(2, 1).IndentLevel = 2
$body.Paragraphs(3, 1).IndentLevel = 2
$body.Paragraphs(4, 1).IndentLevel = 2
$body.Paragraphs(6, 1).IndentLevel = 2

Write-Output "Added slide $($s.SlideIndex) (SlideID=$($s.SlideID)), body paragraphs=$($body.Paragraphs().Count)"
